$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.514.41"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.739.78"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4531"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.076"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.904"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.036"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "1.736.65"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001053"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06347"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "27.542.41"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.087"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "1.938.23"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.034"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.046"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09051"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.646"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.380"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02265"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05947"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2052"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6227"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.880"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.373"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.682"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.697"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5776"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06836"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.78%  "
